$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3088.7778
$ws.Range("I64").Value = 3049.75
$ws.Range("J64").Value = 3120
$ws.Range("K64").Value = 3049.75
$ws.Range("L64").Value = 3120
$ws.Range("M64").Value = -2801.75
$ws.Range("N64").Value = -3616

$ws.Range("H67").Value = 3088.7778
$ws.Range("I67").Value = 3049.75
$ws.Range("J67").Value = 3120
$ws.Range("K67").Value = 3049.75
$ws.Range("L67").Value = 3120
$ws.Range("M67").Value = -2191.75
$ws.Range("N67").Value = -4836

$ws.Range("H70").Value = 45610.06
$ws.Range("I70").Value = 2450
$ws.Range("J70").Value = 51364.734
$ws.Range("K70").Value = 7350
$ws.Range("L70").Value = 154094.202
$ws.Range("M70").Value = -7080
$ws.Range("N70").Value = -154634.202

$ws.Range("H73").Value = 45610.06
$ws.Range("I73").Value = 2450
$ws.Range("J73").Value = 51364.734
$ws.Range("K73").Value = 7350
$ws.Range("L73").Value = 154094.202
$ws.Range("M73").Value = -6414
$ws.Range("N73").Value = -155966.202

$ws.Range("H98").Value = 5607.7
$ws.Range("I98").Value = 5119.6665
$ws.Range("J98").Value = 10000
$ws.Range("K98").Value = 5119.6665
$ws.Range("L98").Value = 10000
$ws.Range("M98").Value = -3621.6665
$ws.Range("N98").Value = -12996

$ws.Range("H103").Value = 396
$ws.Range("I103").Value = 373.5
$ws.Range("J103").Value = 426
$ws.Range("K103").Value = 1120.5
$ws.Range("L103").Value = 1278
$ws.Range("M103").Value = -534.5
$ws.Range("N103").Value = -2450

$ws.Range("H112").Value = 1787.6875
$ws.Range("J112").Value = 1957.4286
$ws.Range("L112").Value = 5872.2858
$ws.Range("N112").Value = -8088.2858

$ws.Range("H116").Value = 4496.5
$ws.Range("I116").Value = 4496.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4496.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1054.5
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 5607.7
$ws.Range("I122").Value = 5119.6665
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 15358.9995
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -12908.9995
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 1452.1154
$ws.Range("I132").Value = 1406.5
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 4219.5
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -1689.5
$ws.Range("N132").Value = -11058.5

$ws.Range("H138").Value = 8288.032
$ws.Range("I138").Value = 5181.727
$ws.Range("J138").Value = 9996.5
$ws.Range("K138").Value = 15545.181
$ws.Range("L138").Value = 29989.5
$ws.Range("M138").Value = -10405.181
$ws.Range("N138").Value = -40269.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1200
$ws.Range("I74").Value = 1200
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1200
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -326
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1200
$ws.Range("I77").Value = 1200
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1632
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H105").Value = 4049.3333
$ws.Range("I105").Value = 4144.727
$ws.Range("K105").Value = 4144.727
$ws.Range("M105").Value = -2397.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 217
$ws.Range("I22").Value = 189.33333
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 189.33333
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 160.66667
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 2145.9768
$ws.Range("I31").Value = 2169.2856
$ws.Range("K31").Value = 2169.2856
$ws.Range("M31").Value = -1874.2856

$ws.Range("H34").Value = 2145.9768
$ws.Range("I34").Value = 2169.2856
$ws.Range("K34").Value = 2169.2856
$ws.Range("M34").Value = -1967.2856

$ws.Range("H60").Value = 15003.363
$ws.Range("I60").Value = 11503.8
$ws.Range("K60").Value = 11503.8
$ws.Range("M60").Value = -10992.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 102.4
$ws.Range("I15").Value = 21.8
$ws.Range("K15").Value = 65.4
$ws.Range("M15").Value = 74.6

$ws.Range("H70").Value = 527
$ws.Range("I70").Value = 527
$ws.Range("K70").Value = 1581
$ws.Range("M70").Value = -1266

$ws.Range("H73").Value = 527
$ws.Range("I73").Value = 527
$ws.Range("K73").Value = 1581
$ws.Range("M73").Value = -489

$ws.Range("H75").Value = 71.75
$ws.Range("I75").Value = 23
$ws.Range("K75").Value = 69
$ws.Range("M75").Value = 929

$ws.Range("H78").Value = 71.75
$ws.Range("I78").Value = 23
$ws.Range("K78").Value = 207
$ws.Range("M78").Value = 4785

$ws.Range("H139").Value = 1764.4
$ws.Range("I139").Value = 1764.4
$ws.Range("K139").Value = 5293.200000000001
$ws.Range("M139").Value = -153.2000000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6328.2856
$ws.Range("I80").Value = 3266
$ws.Range("J80").Value = 8625
$ws.Range("K80").Value = 3266
$ws.Range("L80").Value = 8625
$ws.Range("M80").Value = -2268
$ws.Range("N80").Value = -10621

$ws.Range("H83").Value = 6328.2856
$ws.Range("I83").Value = 3266
$ws.Range("J83").Value = 8625
$ws.Range("K83").Value = 16330
$ws.Range("L83").Value = 43125
$ws.Range("M83").Value = -11338
$ws.Range("N83").Value = -53109

$ws.Range("H122").Value = 33235.562
$ws.Range("I122").Value = 1960
$ws.Range("K122").Value = 5880
$ws.Range("M122").Value = -3430

$ws.Range("H126").Value = 4504
$ws.Range("I126").Value = 4753.5
$ws.Range("J126").Value = 4171.3335
$ws.Range("K126").Value = 14260.5
$ws.Range("L126").Value = 12514.0005
$ws.Range("M126").Value = -11790.5
$ws.Range("N126").Value = -17454.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3183.923
$ws.Range("I7").Value = 3222.2727
$ws.Range("K7").Value = 3222.2727
$ws.Range("M7").Value = -3110.2727

$ws.Range("H16").Value = 14949.5
$ws.Range("J16").Value = 13916.167
$ws.Range("L16").Value = 13916.167
$ws.Range("N16").Value = -14256.167

$ws.Range("H46").Value = 3024.0688
$ws.Range("I46").Value = 2170.3076
$ws.Range("J46").Value = 3717.75
$ws.Range("K46").Value = 2170.3076
$ws.Range("L46").Value = 3717.75
$ws.Range("M46").Value = -1982.3076
$ws.Range("N46").Value = -4093.75

$ws.Range("H126").Value = 3183.923
$ws.Range("I126").Value = 3222.2727
$ws.Range("K126").Value = 9666.8181
$ws.Range("M126").Value = -7196.8181

$ws.Range("H132").Value = 5642.0713
$ws.Range("I132").Value = 4784.857
$ws.Range("J132").Value = 6499.2856
$ws.Range("K132").Value = 14354.571
$ws.Range("L132").Value = 19497.8568
$ws.Range("M132").Value = -11824.571
$ws.Range("N132").Value = -24557.8568

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1987.6666
$ws.Range("I126").Value = 1987.6666
$ws.Range("K126").Value = 5962.9998
$ws.Range("M126").Value = -3492.9998

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H136").Value = 1127.7878
$ws.Range("I136").Value = 1144.1562
$ws.Range("K136").Value = 3432.4686
$ws.Range("M136").Value = -882.4685999999997
